# botek.xlsx fixture update:
#  - "fields" sheet datatype column (C): rename "float" -> "float_qty",
#    "file" -> "image", and the one "int" datatype that carries a unit
#    (crosshole_overlap, %) -> "int_qty".
#  - Move the active selection to C38 (was F1:G1048576).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fields")

# Rows 2..41 hold data; column C is "datatype", column D is "unit".
# New shared strings must be introduced in this order: "image", "float_qty",
# "int_qty" - so pass over the rows once per datatype to match.
$lastRow = 41

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Value2 -eq "file") {
        $ws.Cells.Item($r, 3).Value2 = "image"
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Value2 -eq "float") {
        $ws.Cells.Item($r, 3).Value2 = "float_qty"
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $datatype = $ws.Cells.Item($r, 3).Value2
    $unit = $ws.Cells.Item($r, 4).Value2
    if (($datatype -eq "int") -and ($unit -ne $null) -and ($unit -ne "")) {
        $ws.Cells.Item($r, 3).Value2 = "int_qty"
    }
}

$ws.Range("C38").Select()
